$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 192-201 were re-shuffled (match data re-ordered while the Indice/
# pais/torneio/temporada/data_partida columns A:E stay put). Snapshot the
# "home..url" columns (F:V) for every affected row first, then write them
# back out in their new positions, so the in-place rewrites can't clobber
# a value we still need to read later in the script.
# ---------------------------------------------------------------------------
$snap192 = $ws.Range("F192:V192").Value2
$snap193 = $ws.Range("F193:V193").Value2
$snap195 = $ws.Range("F195:V195").Value2
$snap196 = $ws.Range("F196:V196").Value2
$snap197 = $ws.Range("F197:V197").Value2
$snap198 = $ws.Range("F198:V198").Value2
$snap199 = $ws.Range("F199:V199").Value2
$snap200 = $ws.Range("F200:V200").Value2
$snap201 = $ws.Range("F201:V201").Value2

$ws.Range("F192:V192").Value = $snap195
$ws.Range("F193:V193").Value = $snap196
$ws.Range("F195:V195").Value = $snap193
$ws.Range("F196:V196").Value = $snap192
$ws.Range("F197:V197").Value = $snap201
$ws.Range("F198:V198").Value = $snap200
$ws.Range("F199:V199").Value = $snap197
$ws.Range("F200:V200").Value = $snap198
$ws.Range("F201:V201").Value = $snap199

# ---------------------------------------------------------------------------
# Two brand-new matches were appended at the end of the table: rows 212 and
# 213 (Indice 211 and 212). Clone the formatting of the last existing row
# (211) so the new rows inherit the same styles (bold/bordered index column,
# date-formatted data_partida column, etc.), then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A211:V211").Copy()
$ws.Range("A212:V212").PasteSpecial(-4122)
$ws.Range("A211:V211").Copy()
$ws.Range("A213:V213").PasteSpecial(-4122)

# pais / torneio / temporada repeat verbatim on every row - copy their exact
# (text-typed) values across instead of retyping "2023" and risking it being
# reinterpreted as a number.
$ws.Range("B211:D211").Copy()
$ws.Range("B212:D212").PasteSpecial(-4163)
$ws.Range("B211:D211").Copy()
$ws.Range("B213:D213").PasteSpecial(-4163)

# Row 212: Millonarios 1 x 0 Ind. Medellin, played 2023-11-26
$ws.Cells.Item(212, 1).Value = 211
$ws.Cells.Item(212, 5).Value = 45256.91666666666
$ws.Cells.Item(212, 6).Value = "Millonarios"
$ws.Cells.Item(212, 7).Value = 1
$ws.Cells.Item(212, 8).Value = "Ind. Medellin"
$ws.Cells.Item(212, 9).Value = 0
$ws.Cells.Item(212, 10).Value = 1.99
$ws.Cells.Item(212, 11).Value = "20/11/2023 00:42"
$ws.Cells.Item(212, 12).Value = 1.93
$ws.Cells.Item(212, 13).Value = "26/11/2023 21:55"
$ws.Cells.Item(212, 14).Value = 3.31
$ws.Cells.Item(212, 15).Value = "20/11/2023 00:42"
$ws.Cells.Item(212, 16).Value = 3.36
$ws.Cells.Item(212, 17).Value = "26/11/2023 21:55"
$ws.Cells.Item(212, 18).Value = 4.18
$ws.Cells.Item(212, 19).Value = "20/11/2023 00:42"
$ws.Cells.Item(212, 20).Value = 4.53
$ws.Cells.Item(212, 21).Value = "26/11/2023 21:55"
$ws.Cells.Item(212, 22).Value = "https://www.betexplorer.com/football/colombia/primera-a/millonarios-ind-medellin/0jZMYZw5/"

# Row 213: Atl. Nacional 1 x 0 America De Cali, played 2023-11-27
$ws.Cells.Item(213, 1).Value = 212
$ws.Cells.Item(213, 5).Value = 45257.02083333334
$ws.Cells.Item(213, 6).Value = "Atl. Nacional"
$ws.Cells.Item(213, 7).Value = 1
$ws.Cells.Item(213, 8).Value = "America De Cali"
$ws.Cells.Item(213, 9).Value = 0
$ws.Cells.Item(213, 10).Value = 2.27
$ws.Cells.Item(213, 11).Value = "20/11/2023 00:42"
$ws.Cells.Item(213, 12).Value = 2.39
$ws.Cells.Item(213, 13).Value = "27/11/2023 00:24"
$ws.Cells.Item(213, 14).Value = 3.22
$ws.Cells.Item(213, 15).Value = "20/11/2023 00:42"
$ws.Cells.Item(213, 16).Value = 3.33
$ws.Cells.Item(213, 17).Value = "27/11/2023 00:21"
$ws.Cells.Item(213, 18).Value = 3.41
$ws.Cells.Item(213, 19).Value = "20/11/2023 00:42"
$ws.Cells.Item(213, 20).Value = 3.19
$ws.Cells.Item(213, 21).Value = "27/11/2023 00:21"
$ws.Cells.Item(213, 22).Value = "https://www.betexplorer.com/football/colombia/primera-a/atl-nacional-america-de-cali/EgwRXggB/"
